$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text (e.g. "60.316.12" using
# "." as a thousands separator). Force each target cell's NumberFormat to
# Text first so Excel does not reinterpret numeric-looking values (e.g.
# "0.420") as floating point numbers and silently drop significant
# trailing zeros / precision.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.316.12"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "2.590.92"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "509.65"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").Value = "153.87"
$ws.Range("E6").Value = "  -3.33%  "

$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -5.27%  "

$ws.Range("D9").Value = "2.600.00"
$ws.Range("E9").Value = "  -2.84%  "

$ws.Range("D10").Value = "6.71"

$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("E12").Value = "  -1.38%  "

$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").Value = "3.049.05"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").Value = "60.321.24"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").Value = "21.62"
$ws.Range("E16").Value = "  -3.27%  "

$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "2.596.73"
$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").Value = "350.85"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").Value = "10.53"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "60.20"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "19.41"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").Value = "152.09"
$ws.Range("E32").Value = "  -3.44%  "

$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").Value = "4.02"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("E36").Value = "  -3.75%  "

$ws.Range("D37").Value = "0.860"
$ws.Range("E37").Value = "  +2.56%  "

$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  -4.19%  "

$ws.Range("D39").Value = "0.844"
$ws.Range("E39").Value = "  -4.49%  "

$ws.Range("D40").Value = "36.23"
$ws.Range("E40").Value = "  +1.60%  "

$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  -1.38%  "

$ws.Range("D42").Value = "300.06"
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").Value = "0.620"
$ws.Range("E44").Value = "  -4.33%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0556"
$ws.Range("E46").Value = "  -3.97%  "

$ws.Range("D47").Value = "19.78"

$ws.Range("D48").Value = "4.83"
$ws.Range("E48").Value = "  -4.15%  "

$ws.Range("D49").Value = "0.0232"
$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("D50").Value = "10.30"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").Value = "1.995.18"
$ws.Range("E51").Value = "  -1.84%  "
